# Update DASHBOARD_main/data_kendaraan.xlsx (Sheet1) so that it shows a
# single consolidated vehicle/tax record (profil + kendaraan) instead of the
# previous multi-row tax list, and add the new vehicle-detail columns
# (Nomor_Rangka, Merek, Model, Warna) ahead of the Status column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-obsolete rows 3-8 -----------------------------------
# Only one data row survives (it becomes row 2), so remove the rest of the
# old table body.
$ws.Range("A3:A8").EntireRow.Delete()

# --- Make room for the 4 new vehicle-detail columns --------------------
# They land between the existing "Pajak" (G) / "Status" (H) columns, so
# insert 4 blank columns at I:L; the old H column (and its data) stays put
# and gets repurposed below.
$ws.Range("I1:L1").EntireColumn.Insert()

# --- Header row (row 1) -------------------------------------------------
$ws.Range("H1").Value = "Nomor_Rangka"
$ws.Range("I1").Value = "Merek"
$ws.Range("J1").Value = "Model"
$ws.Range("K1").Value = "Warna"
$ws.Range("L1").Value = "Status"

# --- Data row (row 2) ----------------------------------------------------
# NIK/Plat/Nomor Rangka etc look numeric, so force text formatting first
# (then reset the style back to Normal) to avoid Excel auto-converting
# them to numbers while still keeping them stored as text, matching the
# source data which stores every cell as text.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1234456278949542"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "BG4576HI"
$ws.Range("C2").Value = "Nia Rahmadani"
$ws.Range("D2").Value = "Palembang"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "65000"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "2026-08-02 00:00:00"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "65000"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = "NMR123XYZ"
$ws.Range("I2").Value = "Honda"
$ws.Range("J2").Value = "Sepeda Motor"
$ws.Range("K2").Value = "Biru"
$ws.Range("L2").Value = "LUNAS"
